$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.207.59"
$ws.Range("E2").Value = "  -4.79%  "
$ws.Range("D3").Value = "2.232.41"
$ws.Range("E3").Value = "  -5.75%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.13%  "
$ws.Range("E7").Value = "  -6.71%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.562"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("E12").Value = "  -9.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.11%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "2.571.55"
$ws.Range("E15").Value = "  -5.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.99%  "
$ws.Range("E17").Value = "  -6.26%  "
$ws.Range("D18").Value = "2.231.11"
$ws.Range("E18").Value = "  -5.72%  "
$ws.Range("D19").Value = "43.046.28"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -9.37%  "
$ws.Range("E22").Value = "  -10.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.01%  "
$ws.Range("E24").Value = "  -10.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "238.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.04%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  -9.24%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.31%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0873"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.53%  "
$ws.Range("E37").Value = "  +9.03%  "
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("E39").Value = "  -6.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("E41").Value = "  -11.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0324"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "1.799.89"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.57%  "
$ws.Range("E48").Value = "  -8.97%  "
$ws.Range("E49").Value = "  -7.35%  "
$ws.Range("E50").Value = "  -8.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.20%  "
